$wb = $excel.ActiveWorkbook

$wsUEC = $wb.Worksheets.Item("Unique_EventCodes")
$wsUEC.Columns.Item(1).ColumnWidth = 21.166666666666668

$wsUEC.Range("A2").Value = 'AMH-REP-00000'
$wsUEC.Range("C2").Value = 'INFO'
$wsUEC.Range("D2").Value = 'Report with definition ''Non-Final-Transactions'' is executed by user: Everflow'
$wsUEC.Range("A3").Value = 'AMH-SVC-000160'
$wsUEC.Range("C3").Value = 'WARN'
$wsUEC.Range("D3").Value = 'Alerted service for SnF Queue rbosgb21_finplussrp1ll!pip:0923045, using Output Channel Dnbsop21 [p:0923045].'
$wsUEC.Range("A4").Value = 'AMH-SWF-00501'
$wsUEC.Range("C4").Value = 'ERROR'
$wsUEC.Range("D4").Value = 'SnF session lost for SnF Queue rbosgb21_finplussrp1ll!pip:0923047, SAG Connection:SA_SAG3_LNI'
$wsUEC.Range("A5").Value = 'AMH-SWF-00520'
$wsUEC.Range("C5").Value = 'ERROR'
$wsUEC.Range("D5").Value = 'SnF session lost for SnF Queue rbosgb21_finpllussrp1ll!pip:0923045, using Output Channel Dnbsop21 [p:0923045], reason:SnF Session was closed. Severity: Logic. Parameters [Parameters: [Parameter: [statusAttributes[0] code=Sw.Snf.InvalidSessionId, text=The session identifier is not valid or does not match the current session identifier on the queue, the session has been closed. severity=logic, parameters=[Parameters [Parameter: [Parameter: rbosgb21_finplussrpl1ll!pip:0923251], action=Check that the session identifier is correct or reacquire the queue. ] SAG Connection:SA_SAG1_LNI'
$wsUEC.Range("A6").Value = 'AMH-USERACTION-00027'
$wsUEC.Range("C6").Value = 'INFO'
$wsUEC.Range("D6").Value = '- service-2025-10-28T00:00:00.000Z  The System Job ''Zenabup-completed-workflow'' was executed by user service'
$wsUEC.Range("A7").Value = 'No_AMH_Log_Code'
$wsUEC.Range("C7").Value = 'INFO'
$wsUEC.Range("D7").Value = '- service-2025-10-28T00:00:00.000Z MB0749 Housekeeping task [cleanup-completed-workflow] started.'

$wsLDS = $wb.Worksheets.Item("Log_Details_From_SWIFT")
$wsLDS.Range("A2").Value = 'AMH-REP-00000'
$wsLDS.Range("B2").Value = 'Data not available in Official SWIFT Log Guide'
$wsLDS.Range("C2").Value = 'Data not available in Official SWIFT Log Guide'
$wsLDS.Range("A3").Value = 'AMH-SVC-000160'
$wsLDS.Range("B3").Value = 'Alerted service {serviceCode} on LN {logicalNode} and PN {physicalNode}: {message}'
$wsLDS.Range("C3").Value = 'service in alerted state.'
$wsLDS.Range("A4").Value = 'AMH-SWF-00501'
$wsLDS.Range("B4").Value = 'Data not available in Official SWIFT Log Guide'
$wsLDS.Range("C4").Value = 'Data not available in Official SWIFT Log Guide'
$wsLDS.Range("A5").Value = 'AMH-SWF-00520'
$wsLDS.Range("B5").Value = 'SnF session lost for {object}: {sessionId}, reason:{reason}, SAG Connection:{sagConnection}'
$wsLDS.Range("C5").Value = 'SnF session stopped'
$wsLDS.Range("A6").Value = 'AMH-USERACTION-00027'
$wsLDS.Range("B6").Value = 'The {EntityName} {Code} was executed by user {User}'
$wsLDS.Range("C6").Value = 'A user executed entity code'
$wsLDS.Range("A7").Value = 'No_AMH_Log_Code'
$wsLDS.Range("B7").Value = 'Data not available in Official SWIFT Log Guide'
$wsLDS.Range("C7").Value = 'Data not available in Official SWIFT Log Guide'

$wsMS = $wb.Worksheets.Item("Merged_Summary")
$wsMS.Columns.Item(1).ColumnWidth = 21.166666666666668
$wsMS.Columns.Item(7).ColumnWidth = 92.16666666666667

for ($i = 0; $i -lt 5; $i++) {
    $wsMS.Rows.Item(2).Copy()
    $wsMS.Rows.Item(3).Insert()
}

$wsMS.Range("A2").Value = 'AMH-REP-00000'
$wsMS.Range("B2").Value = 'INFO'
$wsMS.Range("G2").Value = 'Data not available in Official SWIFT Log Guide'
$wsMS.Range("H2").Value = 'Data not available in Official SWIFT Log Guide'
$wsMS.Range("A3").Value = 'AMH-SVC-000160'
$wsMS.Range("B3").Value = 'WARN'
$wsMS.Range("G3").Value = 'Alerted service {serviceCode} on LN {logicalNode} and PN {physicalNode}: {message}'
$wsMS.Range("H3").Value = 'service in alerted state.'
$wsMS.Range("A4").Value = 'AMH-SWF-00501'
$wsMS.Range("B4").Value = 'ERROR'
$wsMS.Range("G4").Value = 'Data not available in Official SWIFT Log Guide'
$wsMS.Range("H4").Value = 'Data not available in Official SWIFT Log Guide'
$wsMS.Range("A5").Value = 'AMH-SWF-00520'
$wsMS.Range("B5").Value = 'ERROR'
$wsMS.Range("G5").Value = 'SnF session lost for {object}: {sessionId}, reason:{reason}, SAG Connection:{sagConnection}'
$wsMS.Range("H5").Value = 'SnF session stopped'
$wsMS.Range("A6").Value = 'AMH-USERACTION-00027'
$wsMS.Range("B6").Value = 'INFO'
$wsMS.Range("G6").Value = 'The {EntityName} {Code} was executed by user {User}'
$wsMS.Range("H6").Value = 'A user executed entity code'
$wsMS.Range("A7").Value = 'No_AMH_Log_Code'
$wsMS.Range("B7").Value = 'INFO'
$wsMS.Range("G7").Value = 'Data not available in Official SWIFT Log Guide'
$wsMS.Range("H7").Value = 'Data not available in Official SWIFT Log Guide'

$wsMS.Range("A1:H7").AutoFilter() | Out-Null
$wsMS.Range("A1:H7").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Merged_Summary!_FilterDatabase") {
        $n.RefersTo = "='Merged_Summary'!`$A`$1:`$H`$7"
    }
}

